$d = $word.ActiveDocument

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "GIS & Geospatial Analysis Consulting") {
        $targetParagraph = $p
        break
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not find paragraph 'GIS & Geospatial Analysis Consulting'"
}

$newText = "• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels`r• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide`r• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis`r"

$r = $targetParagraph.Range
$r.Collapse(0)
$r.InsertAfter($newText)

Write-Output "Inserted three bullet paragraphs after 'GIS & Geospatial Analysis Consulting'"
